$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '23.672.34'
$ws.Range("E2").Value = '  +0.96%  '

$ws.Range("D3").Value = '1.649.88'
$ws.Range("E3").Value = '  +1.22%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9941'

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.9954'
$ws.Range("E5").Value = '  -0.57%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '304.51'
$ws.Range("E6").Value = '  +0.08%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3805'
$ws.Range("E7").Value = '  +0.76%  '

$ws.Range("E8").Value = '  +0.44%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3626'
$ws.Range("E9").Value = '  -0.91%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.245'
$ws.Range("E10").Value = '  +1.51%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08210'
$ws.Range("E11").Value = '  -0.33%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.9954'
$ws.Range("E12").Value = '  -0.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '22.55'
$ws.Range("E13").Value = '  +0.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.554'
$ws.Range("E14").Value = '  -0.05%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.379'
$ws.Range("E15").Value = '  +1.40%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.00001233'
$ws.Range("E16").Value = '  -1.48%  '

$ws.Range("D17").Value = '1.644.32'
$ws.Range("E17").Value = '  +0.99%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '96.86'
$ws.Range("E18").Value = '  +2.90%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06970'
$ws.Range("E19").Value = '  -0.22%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '6.763'
$ws.Range("E20").Value = '  +4.74%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.65'
$ws.Range("E21").Value = '  -0.51%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.9946'
$ws.Range("E22").Value = '  -0.67%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '12.61'
$ws.Range("E23").Value = '  -0.80%  '

$ws.Range("D24").Value = '23.660.80'
$ws.Range("E24").Value = '  +0.90%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.518'
$ws.Range("E25").Value = '  +1.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.139'
$ws.Range("E26").Value = '  -2.20%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.39'
$ws.Range("E27").Value = '  +0.11%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '152.40'
$ws.Range("E28").Value = '  +1.60%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '5.191'
$ws.Range("E29").Value = '  -2.22%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '134.71'
$ws.Range("E30").Value = '  +0.63%  '

$ws.Range("D31").Value = '1.827.66'
$ws.Range("E31").Value = '  +0.85%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.775'
$ws.Range("E32").Value = '  -0.36%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.099'
$ws.Range("E33").Value = '  +7.73%  '

$ws.Range("B34").Value = 'WEMIXTOKEN'
$ws.Range("C34").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '2.037'
$ws.Range("E34").Value = '  -10.43%  '

$ws.Range("B35").Value = 'FraxShare'
$ws.Range("C35").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '11.43'
$ws.Range("E35").Value = '  +6.11%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.02801'
$ws.Range("E36").Value = '  +0.46%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.2519'
$ws.Range("E37").Value = '  -0.22%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.08807'
$ws.Range("E38").Value = '  +0.67%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '6.093'
$ws.Range("E39").Value = '  +1.12%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.07029'
$ws.Range("E40").Value = '  -1.25%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '12.79'
$ws.Range("E41").Value = '  +4.45%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.7080'
$ws.Range("E42").Value = '  +0.36%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.324'
$ws.Range("E43").Value = '  -2.18%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '15.78'
$ws.Range("E44").Value = '  -3.37%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6531'
$ws.Range("E45").Value = '  -0.34%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.348'
$ws.Range("E46").Value = '  +1.03%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9944'
$ws.Range("E47").Value = '  -0.61%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '3.969'
$ws.Range("E48").Value = '  -0.33%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.07987'
$ws.Range("E49").Value = '  -0.44%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '128.05'
$ws.Range("E50").Value = '  +1.72%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.188'
$ws.Range("E51").Value = '  -1.01%  '

